# Admin orderplacements - testcases 2, Testdata, XML File, Helper - Commit
#
# Target sheet: "Mogento,orderplacement" (sheet7.xml)
# Insert 6 new columns (K:P) for OXO product test data, populate the new
# header row (row 1) and a brand new data row (row 8), then fix up the
# selection/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mogento,orderplacement")
$ws.Activate()

# --- Insert 6 new columns before column K, shifting K:AL -> Q:AR ---
$ws.Columns("K:P").Insert()

# The engine's whole-column Insert() sometimes materializes phantom blank
# cells (with inherited column style) past the end of a sparse row's real
# data when the row has a populated cell immediately to the left of the
# insertion point. Rows 6 and 7 hit that case here; a real Excel insert
# would not create those cells, so scrub them back out.
$ws.Range("K6:P7").Clear()

# --- New header values + new data row 8 (OXO product test data) ---
# Order matches original authoring order so shared-string table indices
# line up with the source edit.
$ws.Range("A8").Value = "OXOProducts"
$ws.Range("K1").Value = "SKUNumberoxosimple"
$ws.Range("P1").Value = "oxoQuantity"
$ws.Range("M1").Value = "SKUNumberoxoconfigurable"
$ws.Range("L1").Value = "SKUNumberoxobundle"
$ws.Range("K8").Value = "'32480"
$ws.Range("M8").Value = "61132200B"
$ws.Range("L8").Value = "oxo-13-piece-grilling-set"
$ws.Range("O1").Value = "SKUnumb"
$ws.Range("N1").Value = "Choosecolor"
$ws.Range("N8").Value = "Tot Navy"
$ws.Range("O8").Value = 61132300
$ws.Range("P8").Value = "'3"

# --- Restore/update the sheet's view state (selection) ---
$ws.Range("M12").Select()
